$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D keeps text storage (values look numeric) before assigning new values
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '27.577.97'
$ws.Range("E2").Value = '  -2.22%  '
$ws.Range("D3").Value = '1.849.49'
$ws.Range("E3").Value = '  -3.04%  '
$ws.Range("D4").Value = '1.004'
$ws.Range("E4").Value = '  -0.84%  '
$ws.Range("D5").Value = '333.42'
$ws.Range("E5").Value = '  -1.30%  '
$ws.Range("D6").Value = '1.004'
$ws.Range("E6").Value = '  -0.84%  '
$ws.Range("D7").Value = '0.4615'
$ws.Range("E7").Value = '  -4.93%  '
$ws.Range("E8").Value = '  -3.99%  '
$ws.Range("D9").Value = '46.44'
$ws.Range("E9").Value = '  -2.56%  '
$ws.Range("D10").Value = '0.07922'
$ws.Range("E10").Value = '  -2.26%  '
$ws.Range("D11").Value = '0.9889'
$ws.Range("E11").Value = '  -4.27%  '
$ws.Range("D12").Value = '21.36'
$ws.Range("E12").Value = '  -3.56%  '
$ws.Range("D13").Value = '1.854.39'
$ws.Range("E13").Value = '  -2.79%  '
$ws.Range("D14").Value = '5.934'
$ws.Range("E14").Value = '  -1.97%  '
$ws.Range("D15").Value = '7.084'
$ws.Range("E15").Value = '  -2.92%  '
$ws.Range("D16").Value = '1.006'
$ws.Range("E16").Value = '  -0.80%  '
$ws.Range("E17").Value = '  -0.52%  '
$ws.Range("D18").Value = '0.06646'
$ws.Range("E18").Value = '  -2.36%  '
$ws.Range("D19").Value = '0.00001032'
$ws.Range("E19").Value = '  -2.33%  '
$ws.Range("D20").Value = '17.02'
$ws.Range("E20").Value = '  -1.33%  '
$ws.Range("D21").Value = '1.004'
$ws.Range("E21").Value = '  -0.79%  '
$ws.Range("D22").Value = '27.589.87'
$ws.Range("E22").Value = '  -2.04%  '
$ws.Range("D23").Value = '5.367'
$ws.Range("E23").Value = '  -3.47%  '
$ws.Range("D24").Value = '10.93'
$ws.Range("E24").Value = '  -1.50%  '
$ws.Range("D25").Value = '2.299'
$ws.Range("E25").Value = '  -2.71%  '
$ws.Range("D26").Value = '157.66'
$ws.Range("E26").Value = '  -1.80%  '
$ws.Range("D27").Value = '19.48'
$ws.Range("E27").Value = '  -3.76%  '
$ws.Range("D28").Value = '2.080'
$ws.Range("E28").Value = '  -2.74%  '
$ws.Range("D29").Value = '5.332'
$ws.Range("E29").Value = '  -5.12%  '
$ws.Range("D30").Value = '119.68'
$ws.Range("E30").Value = '  -2.54%  '
$ws.Range("D31").Value = '0.9626'
$ws.Range("E31").Value = '  -2.45%  '
$ws.Range("D32").Value = '0.09331'
$ws.Range("E32").Value = '  -3.51%  '
$ws.Range("D33").Value = '3.552'
$ws.Range("E33").Value = '  -2.58%  '
$ws.Range("D34").Value = '5.282'
$ws.Range("E34").Value = '  -2.36%  '
$ws.Range("D35").Value = '1.338'
$ws.Range("E35").Value = '  -3.11%  '
$ws.Range("B36").Value = 'VeChain'
$ws.Range("C36").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D36").Value = '0.02208'
$ws.Range("E36").Value = '  -2.89%  '
$ws.Range("B37").Value = 'Hedera'
$ws.Range("C37").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D37").Value = '0.05941'
$ws.Range("E37").Value = '  -3.29%  '
$ws.Range("D38").Value = '8.131'
$ws.Range("E38").Value = '  -1.92%  '
$ws.Range("D39").Value = '1.169'
$ws.Range("E39").Value = '  -3.94%  '
$ws.Range("D40").Value = '0.5867'
$ws.Range("E40").Value = '  -2.73%  '
$ws.Range("D41").Value = '0.1853'
$ws.Range("E41").Value = '  -3.54%  '
$ws.Range("D42").Value = '10.21'
$ws.Range("E42").Value = '  -2.72%  '
$ws.Range("D43").Value = '1.243'
$ws.Range("E43").Value = '  -2.66%  '
$ws.Range("D44").Value = '0.5553'
$ws.Range("E44").Value = '  -3.12%  '
$ws.Range("D45").Value = '12.12'
$ws.Range("E45").Value = '  -1.26%  '
$ws.Range("D46").Value = '1.883'
$ws.Range("E46").Value = '  -3.95%  '
$ws.Range("D47").Value = '0.06670'
$ws.Range("E47").Value = '  -2.86%  '
$ws.Range("D48").Value = '110.41'
$ws.Range("E48").Value = '  -3.61%  '
$ws.Range("D49").Value = '1.050'
$ws.Range("E49").Value = '  -3.12%  '
$ws.Range("D50").Value = '1.004'
$ws.Range("E50").Value = '  -0.96%  '
$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D51").Value = '70.11'
$ws.Range("E51").Value = '  -2.76%  '

# Restore default style on column D so no stray number-format style id is left on cells
$ws.Range("D2:D51").Style = "Normal"
